# Insert a new price-observation row for Femacal de La Calera / Arándano (blue)
# right after the existing row 284 (i.e. it becomes the new row 285), pushing
# the former rows 285-297 down to 286-298.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 285..297 down to 286..298, leaving a blank row at 285.
$ws.Rows.Item(285).Insert()

# Populate the newly inserted row 285 with the new record.
$ws.Cells.Item(285, 1).Value  = 3
$ws.Cells.Item(285, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(285, 3).Value  = "Coquimbo"
$ws.Cells.Item(285, 4).Value  = 44931
$ws.Cells.Item(285, 5).Value  = 5
$ws.Cells.Item(285, 6).Value  = "Fruta"
$ws.Cells.Item(285, 7).Value  = 100101
$ws.Cells.Item(285, 8).Value  = "Berries"
$ws.Cells.Item(285, 9).Value  = 100101001
$ws.Cells.Item(285, 10).Value = "Arándano (blue)"
$ws.Cells.Item(285, 11).Value = "Sin especificar"
$ws.Cells.Item(285, 12).Value = "Primera"
$ws.Cells.Item(285, 13).Value = 56
$ws.Cells.Item(285, 14).Value = 4000
$ws.Cells.Item(285, 15).Value = 4000
$ws.Cells.Item(285, 16).Value = 4000
$ws.Cells.Item(285, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(285, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(285, 19).Value = 2000
$ws.Cells.Item(285, 20).Value = 2
